# ULYSSES-5183 Footnote anchor cannot be set to position 0
#
# Each footnote paragraph currently starts with two tab stops defined
# (positions 0 and 400) and an extra leading "<w:r><w:tab/></w:r>" run
# before the footnote-reference mark. A tab stop at position 0 is
# meaningless (and breaks the anchor), and the leading tab run duplicates
# the tab that already follows the footnote-reference mark. Strip both.
#
# We rebuild each footnote paragraph's content via Range.InsertXML so the
# paragraph properties (tab stops) and runs end up exactly as intended,
# without disturbing the rest of the document.

$d = $word.ActiveDocument

$xmlTemplate = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:tabs><w:tab w:val="start" w:pos="400"/></w:tabs></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:cs="Helvetica" w:eastAsia="Helvetica" w:hAnsi="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:vertAlign w:val="superscript"/></w:rPr><w:footnoteRef/></w:r><w:r><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:cs="Helvetica" w:eastAsia="Helvetica" w:hAnsi="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">{0}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

for ($i = 1; $i -le $d.Footnotes.Count; $i++) {
    $fn = $d.Footnotes($i)
    $rng = $fn.Range

    # The footnote story text begins with: tab, footnote-ref marker, tab,
    # then the visible footnote content - e.g. "\t\x02\tThis is ...".
    $rawText = $rng.Text
    $bodyText = $rawText.Substring(3)

    # Only touch footnotes that actually have the old duplicated-tab shape
    # (defensive - avoids corrupting anything unexpected).
    if ($rawText.Length -ge 3 -and [int][char]$rawText[0] -eq 9 -and [int][char]$rawText[2] -eq 9) {
        $escaped = $bodyText.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
        $xml = $xmlTemplate -f $escaped
        $rng.InsertXML($xml)
    }
}
